$d = $word.ActiveDocument

# 1) Merge the "Mon Sep 17" + " 11:42:44 PDT 2017" runs into a single run.
$null = $d.Content.Find.Execute("Mon Sep 17 11:42:44 PDT 2017", $true, $false, $false, $false, $false, $true, 1, $false, "Mon Sep 17 11:42:44 PDT 2017", 2)

# 2) Insert the new "Tue Sep 18" purchase-details block after the
#    "Amount balance ... - 35430.0" paragraph.
$anchorText = "- 35430.0"
$found = $d.Content.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*- 35430.0*") {
        $anchorPara = $d.Paragraphs($i)
    }
}

$r = $anchorPara.Range
$r.Collapse(0)
$beforeCount = $d.Paragraphs.Count
$r.InsertAfter("`rTue Sep 18 10:56:47 PDT 2017`rPerson Name`t`t`t`t- TRM`r---------------------------------------------------------------`rItem Name`t`t`t`t- CHOW`rNumber of Pockets`t`t`t- 2`rNumber of KGs`t`t`t- 128`rRate`t`t`t`t`t- 15`rTotal Price`t`t`t`t- 1920.0`rAmount balance`t`t`t- 37350.0`r`rItem Name`t`t`t`t- CHOWCHOW`rNumber of Pockets`t`t`t- 1`rNumber of KGs`t`t`t- 62`rRate`t`t`t`t`t- 14`rTotal Price`t`t`t`t- 868.0`rAmount balance`t`t`t- 38218.0`r")
Write-Output "paragraphs before=$beforeCount after=$($d.Paragraphs.Count)"
